{"js": "// Office.js (Word JavaScript API) script.\n// This is the body of `async (context) => { ... }`.\n//\n// Changes applied (per the target diff):\n//  1. Remove the whole \"Due: Thursday June 6 by 23:59pm (emailed)\" paragraph.\n//  2. Collapse the \"Summary: Start with <gram>a brief summary</gram> of the\n//     article...\" run split (and its surrounding grammar-check markers)\n//     into a single contiguous run with the same visible text.\n//  3. Remove the whole \"This assignment is pass/fail based on completion\n//     and adherence to these guidelines.\" paragraph.\n\nconst body = context.document.body;\n\n// --- 1 & 3: delete the two whole paragraphs, found by their leading text ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Due:\") === 0 && text.indexOf(\"23:59pm\") !== -1) {\n    paragraphs.items[i].delete();\n  } else if (text.indexOf(\"This assignment is pass/fail\") === 0) {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n\n// --- 2: normalize the split \"Start with / a brief summary / of the ---\n// --- article...\" run into one run with identical text -----------------\nconst target =\n  \": Start with a brief summary of the article, focusing on its main \" +\n  \"thesis and methodology.\";\nconst hits = body.search(target, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  // Re-insert the identical text over the matched range; Word collapses\n  // the range (which spans the old run boundaries / proofErr marks) into\n  // a single run.\n  hits.items[0].insertText(target, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument / $d are available.\n#\n# Changes applied (per the target diff):\n#  1. Remove the whole \"Due: Thursday June 6 by 23:59pm (emailed)\" paragraph.\n#  2. Collapse the \"Summary: Start with <gram>a brief summary</gram> of the\n#     article...\" run split (and its surrounding grammar-check markers)\n#     into a single contiguous run with the same visible text.\n#  3. Remove the whole \"This assignment is pass/fail based on completion\n#     and adherence to these guidelines.\" paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1: delete the whole \"Due: ...\" paragraph -----------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Due:*23:59pm*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- 3: delete the whole \"This assignment is pass/fail...\" paragraph ------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"This assignment is pass/fail*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- 2: normalize the split \"Start with / a brief summary / of the -------\n# --- article...\" run into one run with identical text ---------------------\n$target = \": Start with a brief summary of the article, focusing on its main thesis and methodology.\"\n$marker = [char]1\n$markedTarget = $target + $marker\n\n# Pass 1: append a one-character marker. Word merges the matched range into\n# a single run as part of the text replacement (the marker forces the\n# engine to recognize the text as changed, instead of a no-op).\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$found1 = $rng1.Find.Execute($target)\nif ($found1) {\n    $rng1.Text = $markedTarget\n}\n\n# Pass 2: find the now-merged run (with trailing marker) and strip the\n# marker back off, leaving the original visible text in a single run.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute($markedTarget)\nif ($found2) {\n    $rng2.Text = $target\n}\n"}
